# Weekly update: insert two new price records (rows 87-88) for
# "Terminal La Palmera de La Serena" - Caqui - Mankaki (Primera/Segunda),
# dated 44711, sourced from Región de O'Higgins. All rows that were
# previously at 87-107 shift down to 89-109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 87; this shifts the
# existing rows 87-107 down to 89-109 and copies formatting (e.g. the
# date style on column D) from the row above, same as Excel's UI does.
$ws.Rows("87:88").Insert()

# ---- Row 87: Mankaki / Primera ----
$ws.Cells.Item(87, 1).Value = 8
$ws.Cells.Item(87, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 44711
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(87, 6).Value = "Fruta"
$ws.Cells.Item(87, 7).Value = 100107
$ws.Cells.Item(87, 8).Value = "Otros"
$ws.Cells.Item(87, 9).Value = 100107001
$ws.Cells.Item(87, 10).Value = "Caqui"
$ws.Cells.Item(87, 11).Value = "Mankaki"
$ws.Cells.Item(87, 12).Value = "Primera"
$ws.Cells.Item(87, 13).Value = 16
$ws.Cells.Item(87, 14).Value = 330000
$ws.Cells.Item(87, 15).Value = 340000
$ws.Cells.Item(87, 16).Value = 335000
$ws.Cells.Item(87, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(87, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(87, 19).Value = 744
$ws.Cells.Item(87, 20).Value = 450

# ---- Row 88: Mankaki / Segunda ----
$ws.Cells.Item(88, 1).Value = 8
$ws.Cells.Item(88, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(88, 3).Value = "Coquimbo"
$ws.Cells.Item(88, 4).Value = 44711
$ws.Cells.Item(88, 5).Value = 4
$ws.Cells.Item(88, 6).Value = "Fruta"
$ws.Cells.Item(88, 7).Value = 100107
$ws.Cells.Item(88, 8).Value = "Otros"
$ws.Cells.Item(88, 9).Value = 100107001
$ws.Cells.Item(88, 10).Value = "Caqui"
$ws.Cells.Item(88, 11).Value = "Mankaki"
$ws.Cells.Item(88, 12).Value = "Segunda"
$ws.Cells.Item(88, 13).Value = 20
$ws.Cells.Item(88, 14).Value = 290000
$ws.Cells.Item(88, 15).Value = 300000
$ws.Cells.Item(88, 16).Value = 295000
$ws.Cells.Item(88, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(88, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(88, 19).Value = 656
$ws.Cells.Item(88, 20).Value = 450
